# Replace the text of a paragraph whose *entire* text matches $oldText with
# $newText, while preserving the paragraph's run structure (including any
# leading empty <w:r/> runs) by round-tripping through WordOpenXML instead
# of using a plain Range.Text / Find-Replace assignment (which this engine,
# like Word itself, normalizes/merges adjacent runs for).
function Replace-ParaText {
    param($doc, [string]$oldText, [string]$newText)

    $updated = 0
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text
        # Paragraph text ends with a paragraph mark (CR) or, for the very
        # last paragraph in a story, a cell/section mark; strip those off
        # before comparing.
        $t = $t.TrimEnd([char]13, [char]7)

        if ($t -eq $oldText) {
            $xml = $p.Range.WordOpenXML

            # WordOpenXML() mints a standalone fragment and stamps the
            # paragraph with fresh w14:paraId/rsid* bookkeeping attributes
            # that are not present in the source document; strip them back
            # out so re-inserting the XML doesn't introduce spurious
            # attributes that weren't part of the intended edit.
            $newXml = $xml -replace [regex]::Escape($oldText), $newText
            $newXml = $newXml -replace ' w14:paraId="[0-9A-Fa-f]+"', ''
            $newXml = $newXml -replace ' w14:textId="[0-9A-Fa-f]+"', ''
            $newXml = $newXml -replace ' w:rsidR="[0-9A-Fa-f]*"', ''
            $newXml = $newXml -replace ' w:rsidRDefault="[0-9A-Fa-f]*"', ''
            $newXml = $newXml -replace ' w:rsidP="[0-9A-Fa-f]*"', ''
            $newXml = $newXml -replace ' w:rsidRPr="[0-9A-Fa-f]*"', ''

            $p.Range.InsertXML($newXml)
            $updated = $updated + 1
        }
    }
    return $updated
}

$d = $word.ActiveDocument

# Heading call-to-action (appears twice: the Heading1 title and the bold
# line near the end of the document).
Replace-ParaText $d "Play Before Time Runs Out for free today!" "Play Before Time Runs Out Free" | Out-Null

# "What we like" bullet list
Replace-ParaText $d "Atmospheric design that immerses players in oriental culture" "Well-designed symbols and background" | Out-Null
Replace-ParaText $d "Unique mechanics of dueling that can transform regular symbols into Wild" "Atmospheric soundtrack" | Out-Null
Replace-ParaText $d "Fixed lines ensure constant chances for winning combinations" "Range of bonus games" | Out-Null
Replace-ParaText $d "Wide range of bonus games and free spins" "Free spins with extra spin awards" | Out-Null

# "What we don't like" bullet list
Replace-ParaText $d "Limited paylines may not appeal to players who prefer more options" "Limited number of paylines" | Out-Null
Replace-ParaText $d "High variance may lead to long stretches of play without significant wins" "High cost per token" | Out-Null

# Closing italic summary line
Replace-ParaText $d "Read our review of Before Time Runs Out, an immersive slot game set in oriental culture, and play for free with exciting bonus games and free spins." "Review of Before Time Runs Out - Play this slot game for free and enjoy bonus games and free spins." | Out-Null
